$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row: Spanish term "Propiedad " / English term "Property"
$ws.Range("A7").Value = "Propiedad "
$ws.Range("B7").Value = "Property"

# English column gets a left-aligned style with one level of indent
$ws.Range("B7").HorizontalAlignment = -4131 # xlLeft
$ws.Range("B7").IndentLevel = 1
